# Generate Report for Handback
# Adds a new handback row for file 733b098d-0dbe-48f9-902d-1e8aa3bd62ea
# (handed back in sync with en-US) to the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$uuid   = "733b098d-0dbe-48f9-902d-1e8aa3bd62ea"
$hash   = "36e5d3ad3d4f3117829f42503aaa090637fd9776"
$mdName = "$uuid.md"
$zhXlf  = "$uuid.$hash.zh-cn.xlf"
$deXlf  = "$uuid.$hash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$includeText  = "Include"

$zhHandoffDt  = "2016-01-27 02:23:38"
$zhHandbackDt = "2016-01-27 02:24:23"
$deHandoffDt  = "2016-01-27 02:23:50"
$deHandbackDt = "2016-01-27 02:24:43"

# ---------------------------------------------------------------
# Sheet "Overview": columns A (File Name, hyperlink), B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------
# Sheet "zh-cn": A Source File Name, B Status, C Correspond Handoff File,
# D Correspond Handoff Datetime, E Target File, F Correspond Handback File,
# G Correspond Handback DateTime, H Handoff Reason, I Dependency From
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Hyperlinks.Add(
    $wsZh.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000001/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf",
    "",
    "",
    $zhXlf
)
$wsZh.Range("D4").Value = $zhHandoffDt
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000002/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000003/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf",
    "",
    "",
    $zhXlf
)
$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $includeText

# ---------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn, for the German locale files
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Hyperlinks.Add(
    $wsDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000004/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf",
    "",
    "",
    $deXlf
)
$wsDe.Range("D4").Value = $deHandoffDt
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000005/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/733b098d0dbe48f9902d1e8aa3bd62ea0000006/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf",
    "",
    "",
    $deXlf
)
$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $includeText

Write-Output "Handback report row added for $uuid"
